# Marksheet update: populate the "Student Ans" columns (A, and D for the
# first few rows) with the student's actual answers, recompute the score
# summary (right / wrong / not-attempted / max / total), and drop the
# now-unused extra "Student Ans / Correct Ans" block that used to live in
# columns G:H (plus the now-empty D:E tail), so the used range shrinks back
# down to A5:E40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the now-unused third answer block (G:H) entirely, and the
#    tail of the second block (D:E) that is no longer needed below row 18.
#    Clear() (not ClearContents) drops the cells completely so the sheet
#    dimension shrinks from A5:H40 down to A5:E40.
# ---------------------------------------------------------------------
$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------
# 2. Summary block (rows 10-12): give the row labels in column A the
#    "mtitleStyle" look (they were unstyled before), and fill in the
#    real right/wrong/not-attempted/max counts plus the marking scheme
#    and the final score.
# ---------------------------------------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

# Row 10: Right / Wrong / Not Attempt / Max
$ws.Range("B10").Value = 24
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 28

# Row 11: Marking scheme (+4 / -1). C11 used to be stored as text "-1";
# store it as a real number now so downstream math doesn't choke on it.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: Total score
$ws.Range("B12").Value = 96
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "95/112"

# ---------------------------------------------------------------------
# 3. Per-question "Student Ans" columns: fill column A (and column D for
#    the rows that still have a second block, 16-18) with the student's
#    actual answer. Where the student answered correctly it mirrors the
#    "Correct Ans" column and is styled "correctStyle"; where they were
#    wrong the wrong answer is shown with "incorrectStyle"; the few
#    questions left unanswered stay blank with the plain "normalStyle".
# ---------------------------------------------------------------------
function Set-StudentAnswer($col, $row, $value, $style) {
    $cell = $ws.Range("$col$row")
    $cell.Value = $value
    $cell.Style = $style
}

# Column A mirrors column B (the "Correct Ans" for the main block).
Set-StudentAnswer "A" 16 "Option A" "correctStyle"
Set-StudentAnswer "A" 17 "Option D" "correctStyle"
Set-StudentAnswer "A" 18 "Option B" "correctStyle"
Set-StudentAnswer "A" 19 "Option C" "correctStyle"
Set-StudentAnswer "A" 20 "Option B" "correctStyle"
Set-StudentAnswer "A" 21 "Option C" "correctStyle"
Set-StudentAnswer "A" 22 "Option D" "correctStyle"
Set-StudentAnswer "A" 23 "Option D" "correctStyle"
Set-StudentAnswer "A" 24 "Option A" "correctStyle"
Set-StudentAnswer "A" 25 "Option B" "incorrectStyle"
Set-StudentAnswer "A" 26 "Option C" "correctStyle"
Set-StudentAnswer "A" 27 "Option A" "correctStyle"
# Row 28 left blank (not attempted) - stays "normalStyle", no change.
Set-StudentAnswer "A" 29 "Option D" "correctStyle"
Set-StudentAnswer "A" 30 "Option B" "correctStyle"
Set-StudentAnswer "A" 31 "Option D" "correctStyle"
Set-StudentAnswer "A" 32 "Option C" "correctStyle"
Set-StudentAnswer "A" 33 "Option D" "correctStyle"
# Row 34 left blank (not attempted) - stays "normalStyle", no change.
Set-StudentAnswer "A" 35 "Option D" "correctStyle"
Set-StudentAnswer "A" 36 "Option A" "correctStyle"
Set-StudentAnswer "A" 37 "Option A" "correctStyle"
Set-StudentAnswer "A" 38 "Option A" "correctStyle"
Set-StudentAnswer "A" 39 "Option D" "correctStyle"
# Row 40 left blank (not attempted) - stays "normalStyle", no change.

# Column D (rows 16-18 only) mirrors column E (the "Correct Ans" for the
# second block) - the student got all three of these right too.
Set-StudentAnswer "D" 16 "Option A" "correctStyle"
Set-StudentAnswer "D" 17 "Option C" "correctStyle"
Set-StudentAnswer "D" 18 "Option D" "correctStyle"
